$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update the test case name in C2 from "testSubscribe" to "testWorkFrontJobs"
$ws.Range("C2").Value = "testWorkFrontJobs"

# Select C2 to match the resulting selection in the sheet view
$ws.Range("C2").Select()
